$d = $word.ActiveDocument

# The stakeholder-analysis table is the 2nd table in the document
# (Stakeholder | Value/Benefit | Interest), which becomes
# (Stakeholder | Value/Benefit | Attitude | Key Interests | Constraints).
$t = $d.Tables.Item(2)

# Helper to set a cell's text while keeping a clean <w:r><w:t> run (no
# stray paragraph-mark rPr), mirroring how the other cells are authored.
function Set-CellText($cell, [string]$text, [bool]$bold) {
    $cell.Range.Text = $text
    $start = $cell.Range.Start
    $r = $d.Range($start, $start + $text.Length)
    if ($bold) {
        $r.Font.Bold = $true
    }
}

# Add two new columns so the table goes from 3 to 5 columns.
$t.Columns.Add() | Out-Null
$t.Columns.Add() | Out-Null

# Rename the existing 3rd column header from "Interest" to "Attitude",
# and fill in the two new header cells.
Set-CellText $t.Cell(1,3) "Attitude" $true
Set-CellText $t.Cell(1,4) "Key Interests" $true
Set-CellText $t.Cell(1,5) "Constraints" $true

# Row 2 - Shop Owner: old col3 text moves (extended) to the new col4;
# col3 gets the new "Attitude" value and col5 the new "Constraints" value.
Set-CellText $t.Cell(2,4) "Revenue reports, System security, Cost reduction." $false
Set-CellText $t.Cell(2,3) "Strong supporter" $false
Set-CellText $t.Cell(2,5) "Budget for new hardware." $false

# Row 3 - Administrator
Set-CellText $t.Cell(3,4) "User Management, Data Integrity, Ease of config." $false
Set-CellText $t.Cell(3,3) "Supportive but wary of complexity" $false
Set-CellText $t.Cell(3,5) "Time needed for initial setup." $false

# Row 4 - Sales Staff
Set-CellText $t.Cell(4,4) "Ease of use, Stability, Speed of checkout." $false
Set-CellText $t.Cell(4,3) "Potential resistance to change" $false
Set-CellText $t.Cell(4,5) "Limited technical skills." $false

# All columns become evenly split: 1728 dxa = 86.4 pt each.
for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $t.Cell(1,$c).Width = 86.4
}
